# Fruta / hortaliza, semanal
# Insert a new weekly record as row 53, shifting the existing rows 53-55
# down to rows 54-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 53 (old rows 53-55 -> 54-56)
$ws.Rows.Item(53).Insert()

# Populate the new row 53 with the new weekly record.
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 44585
$ws.Range("D53").NumberFormat = $ws.Range("D54").NumberFormat
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100101
$ws.Range("H53").Value = "Berries"
$ws.Range("I53").Value = 100101001
$ws.Range("J53").Value = "Arándano (blue)"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Segunda"
$ws.Range("M53").Value = 150
$ws.Range("N53").Value = 3000
$ws.Range("O53").Value = 3000
$ws.Range("P53").Value = 3000
$ws.Range("Q53").Value = "$/bandeja 2 kilos"
$ws.Range("R53").Value = "Provincia de Linares"
$ws.Range("S53").Value = 1500
$ws.Range("T53").Value = 2
